$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.1827
$ws.Range("B4").Value = 4.618700000000002
$ws.Range("C4").Value = -11.0877

$ws.Range("B5").Value = 5.272999999999997

$ws.Range("A7").Value = -21.59470000000001

$ws.Range("B8").Value = 4.923099999999999

$ws.Range("C9").Value = -11.72830000000001

$ws.Range("A16").Value = -21.56980000000002
$ws.Range("B16").Value = 4.944500000000002

$ws.Range("C18").Value = -14.68500000000001

$wb.Save()
